# Append two new rows (196, 197) of ECONOMICS:QAM2 data to the sheet,
# following the existing table's formatting/layout (row 195 is the
# last existing data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 196; Date = 44986.45833333334; Value = 704618000000 },
    @{ Row = 197; Date = 45017.45833333334; Value = 696270000000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prevRow = $row - 1

    # Copy the previous row's full formatting (style, borders, number
    # format, etc.) down into the new row, same as dragging the fill
    # handle / inserting a like-formatted row.
    $src = $ws.Range("A" + $prevRow + ":G" + $prevRow)
    $dst = $ws.Range("A" + $row + ":G" + $row)
    $src.Copy($dst)

    # Now overwrite with this row's actual values.
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:QAM2"
    $ws.Cells.Item($row, 3).Value = $r.Value
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 5).Value = $r.Value
    $ws.Cells.Item($row, 6).Value = $r.Value
    $ws.Cells.Item($row, 7).Value = 0
}
